$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tbl = $ws.ListObjects.Item("Tabela1")

$row = $tbl.ListRows.Add()

$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "Aula 3"
$ws.Range("B4").Value = "Discussão sobre tipos de variáveis e desenho de experimentos."
$ws.Range("C4").Value = 45910
$ws.Range("D4").Value = 0.66666666666666663
$ws.Range("E4").Value = 0.72222222222222221
$ws.Range("F4").Formula = "=HOUR(Tabela1[[#This Row],[Horário de fim da aula]]-Tabela1[[#This Row],[Horário de nício da aula]])+(MINUTE(Tabela1[[#This Row],[Horário de fim da aula]]-Tabela1[[#This Row],[Horário de nício da aula]])/60)"
$ws.Range("G4").Value = "Reconhecendo os tipos de variáveis no exerimento no Excel;`nComo organizar os dados no Excel;`nEntendo sobre variáveis dependentes e independentes;`nEntendendo sobre fatores e níveis de fatores;`nDiscussão sobre desenho de experimento (simples e fatorial) e analises estatísticas (paramétrica, não-paramétrica, análise univariada, bivariada e multivariada)."

$ws.Range("G4").RowHeight = 90

$ws.Range("G5").Select() | Out-Null
